$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.606.17"
$ws.Range("E2").Value = "'  -2.24%  "
$ws.Range("D3").Value = "'1.843.35"
$ws.Range("E3").Value = "'  -1.31%  "
$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("D5").Value = "'314.16"
$ws.Range("E5").Value = "'  -1.53%  "
$ws.Range("E6").Value = "'  -0.22%  "
$ws.Range("D7").Value = "'0.4244"
$ws.Range("E7").Value = "'  -3.43%  "
$ws.Range("D8").Value = "'0.3644"
$ws.Range("E8").Value = "'  -1.27%  "
$ws.Range("D9").Value = "'45.32"
$ws.Range("E9").Value = "'  +0.52%  "
$ws.Range("D10").Value = "'0.07263"
$ws.Range("E10").Value = "'  -3.10%  "
$ws.Range("D11").Value = "'0.8906"
$ws.Range("E11").Value = "'  -4.77%  "
$ws.Range("D12").Value = "'20.60"
$ws.Range("E12").Value = "'  -3.45%  "
$ws.Range("D13").Value = "'1.826.24"
$ws.Range("E13").Value = "'  -2.93%  "
$ws.Range("D14").Value = "'5.360"
$ws.Range("E14").Value = "'  -2.04%  "
$ws.Range("D15").Value = "'6.564"
$ws.Range("E15").Value = "'  -2.03%  "
$ws.Range("D16").Value = "'0.06862"
$ws.Range("E16").Value = "'  -0.55%  "
$ws.Range("E17").Value = "'  -0.21%  "
$ws.Range("D18").Value = "'78.49"
$ws.Range("E18").Value = "'  -4.29%  "
$ws.Range("D19").Value = "'0.000008827"
$ws.Range("E19").Value = "'  -2.30%  "
$ws.Range("E20").Value = "'  -0.24%  "
$ws.Range("D21").Value = "'15.45"
$ws.Range("E21").Value = "'  -2.92%  "
$ws.Range("D22").Value = "'27.589.74"
$ws.Range("E22").Value = "'  -2.24%  "
$ws.Range("D23").Value = "'4.996"
$ws.Range("E23").Value = "'  -2.40%  "
$ws.Range("D24").Value = "'10.57"
$ws.Range("E24").Value = "'  -2.15%  "
$ws.Range("D25").Value = "'2.064.47"
$ws.Range("E25").Value = "'  -3.45%  "
$ws.Range("E26").Value = "'  +0.46%  "
$ws.Range("D27").Value = "'155.21"
$ws.Range("E27").Value = "'  -0.05%  "
$ws.Range("E28").Value = "'  +0.63%  "
$ws.Range("B29").Value = "'InternetComputer(DFINITY)"
$ws.Range("C29").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.232"
$ws.Range("E29").Value = "'  -1.62%  "
$ws.Range("B30").Value = "'BitcoinCash"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'118.26"
$ws.Range("E30").Value = "'  +4.27%  "
$ws.Range("D31").Value = "'1.815"
$ws.Range("E31").Value = "'  +5.14%  "
$ws.Range("D32").Value = "'0.08881"
$ws.Range("E32").Value = "'  -1.62%  "
$ws.Range("D33").Value = "'0.7759"
$ws.Range("E33").Value = "'  -2.46%  "
$ws.Range("D34").Value = "'4.567"
$ws.Range("E34").Value = "'  -5.75%  "
$ws.Range("D35").Value = "'2.952"
$ws.Range("E35").Value = "'  +0.54%  "
$ws.Range("D36").Value = "'1.101"
$ws.Range("E36").Value = "'  -6.24%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "'  -0.26%  "
$ws.Range("D38").Value = "'0.05410"
$ws.Range("E38").Value = "'  -0.67%  "
$ws.Range("D39").Value = "'1.099"
$ws.Range("E39").Value = "'  -2.58%  "
$ws.Range("D40").Value = "'0.01917"
$ws.Range("E40").Value = "'  -2.60%  "
$ws.Range("D41").Value = "'2.768"
$ws.Range("E41").Value = "'  -6.15%  "
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.838"
$ws.Range("E42").Value = "'  -2.98%  "
$ws.Range("B43").Value = "'TheSandbox"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5061"
$ws.Range("E43").Value = "'  -3.86%  "
$ws.Range("D44").Value = "'0.1651"
$ws.Range("E44").Value = "'  -1.80%  "
$ws.Range("D45").Value = "'8.190"
$ws.Range("D46").Value = "'0.06617"
$ws.Range("E46").Value = "'  -1.82%  "
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.31"
$ws.Range("E47").Value = "'  -2.08%  "
$ws.Range("D48").Value = "'105.37"
$ws.Range("E48").Value = "'  -1.55%  "
$ws.Range("B49").Value = "'Decentraland"
$ws.Range("C49").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.4697"
$ws.Range("E49").Value = "'  -3.62%  "
$ws.Range("D50").Value = "'0.9998"
$ws.Range("E50").Value = "'  -0.29%  "
$ws.Range("D51").Value = "'1.624"
